$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Se necesitó investigar este tema para crear las url únicas  para cada torneo, edicion, equipo, jugador. "
$ws.Range("A9").Value = "URL Rewriting en ASP.NET"
$ws.Range("C9").Value = "Antonio"
$ws.Range("D9").Value = "No"
$ws.Range("E9").Value = "#Sprint 13"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "Link a info https://msdn.microsoft.com/en-us/library/ms972974.aspx"
$ws.Range("H9").Value = "Realizada"
